# Scheduled runner update: refresh Atomos Profits price/profit columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 825.05
$ws.Range("I33").Value = 700
$ws.Range("J33").Value = 1116.8334
$ws.Range("K33").Value = 700
$ws.Range("L33").Value = 1116.8334
$ws.Range("M33").Value = -471
$ws.Range("N33").Value = -1574.8334
$ws.Range("H111").Value = 1077.8889
$ws.Range("I111").Value = 1075.125
$ws.Range("K111").Value = 3225.375
$ws.Range("M111").Value = -158.375
$ws.Range("H129").Value = 4630677.5
$ws.Range("J129").Value = 1004
$ws.Range("L129").Value = 3012
$ws.Range("N129").Value = -13012
$ws.Range("H131").Value = 2828.423
$ws.Range("I131").Value = 2249.1667
$ws.Range("J131").Value = 4131.75
$ws.Range("K131").Value = 6747.500100000001
$ws.Range("L131").Value = 12395.25
$ws.Range("M131").Value = -1707.500100000001
$ws.Range("N131").Value = -22475.25
$ws.Range("H138").Value = 3587.3613
$ws.Range("I138").Value = 1858.5918
$ws.Range("J138").Value = 6078.8237
$ws.Range("K138").Value = 5575.7754
$ws.Range("L138").Value = 18236.4711
$ws.Range("M138").Value = -435.7753999999995
$ws.Range("N138").Value = -28516.4711

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14343.237
$ws.Range("I32").Value = 13508.792
$ws.Range("K32").Value = 13508.792
$ws.Range("M32").Value = -13221.792
$ws.Range("H97").Value = 446.92307
$ws.Range("I97").Value = 469.16666
$ws.Range("K97").Value = 469.16666
$ws.Range("M97").Value = 26.83334000000002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 26404.8
$ws.Range("J10").Value = 32506
$ws.Range("L10").Value = 32506
$ws.Range("N10").Value = -32786
$ws.Range("H102").Value = 28400
$ws.Range("I102").Value = 20000
$ws.Range("J102").Value = 29800
$ws.Range("K102").Value = 20000
$ws.Range("L102").Value = 29800
$ws.Range("M102").Value = -16755
$ws.Range("N102").Value = -36290
$ws.Range("H141").Value = 33830
$ws.Range("J141").Value = 28440
$ws.Range("L141").Value = 28440
$ws.Range("N141").Value = -38800

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 951.44446
$ws.Range("I22").Value = 344.7143
$ws.Range("K22").Value = 344.7143
$ws.Range("M22").Value = 5.28570000000002
$ws.Range("H31").Value = 4354.283
$ws.Range("I31").Value = 3369.2273
$ws.Range("J31").Value = 5053.355
$ws.Range("K31").Value = 3369.2273
$ws.Range("L31").Value = 5053.355
$ws.Range("M31").Value = -3074.2273
$ws.Range("N31").Value = -5643.355
$ws.Range("H34").Value = 4354.283
$ws.Range("I34").Value = 3369.2273
$ws.Range("J34").Value = 5053.355
$ws.Range("K34").Value = 3369.2273
$ws.Range("L34").Value = 5053.355
$ws.Range("M34").Value = -3167.2273
$ws.Range("N34").Value = -5457.355
$ws.Range("H107").Value = 1374.9688
$ws.Range("I107").Value = 1333.2858
$ws.Range("K107").Value = 1333.2858
$ws.Range("M107").Value = 586.7141999999999
$ws.Range("H141").Value = 28325
$ws.Range("J141").Value = 28325
$ws.Range("L141").Value = 28325
$ws.Range("N141").Value = -38685

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 432.27274
$ws.Range("I12").Value = 25.166666
$ws.Range("J12").Value = 664.9048
$ws.Range("K12").Value = 75.49999800000001
$ws.Range("L12").Value = 1994.7144
$ws.Range("M12").Value = 97.50000199999999
$ws.Range("N12").Value = -2340.7144
$ws.Range("H131").Value = 1575.7142
$ws.Range("J131").Value = 1154.6571
$ws.Range("L131").Value = 3463.9713
$ws.Range("N131").Value = -13543.9713

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5263.3887
$ws.Range("I122").Value = 4199.5864
$ws.Range("J122").Value = 6497.4
$ws.Range("K122").Value = 12598.7592
$ws.Range("L122").Value = 19492.2
$ws.Range("M122").Value = -10148.7592
$ws.Range("N122").Value = -24392.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2914.2856
$ws.Range("I7").Value = 1550
$ws.Range("J7").Value = 3460
$ws.Range("K7").Value = 1550
$ws.Range("L7").Value = 3460
$ws.Range("M7").Value = -1438
$ws.Range("N7").Value = -3684
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H100").Value = 2232.7273
$ws.Range("I100").Value = 1228.5714
$ws.Range("J100").Value = 3990
$ws.Range("K100").Value = 1228.5714
$ws.Range("L100").Value = 3990
$ws.Range("M100").Value = -687.5714
$ws.Range("N100").Value = -5072
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988
$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524
$ws.Range("H121").Value = 40420
$ws.Range("J121").Value = 40420
$ws.Range("L121").Value = 40420
$ws.Range("M121").Value = -43914
$ws.Range("H122").Value = 3403.4333
$ws.Range("I122").Value = 2709.4546
$ws.Range("J122").Value = 5311.875
$ws.Range("K122").Value = 8128.3638
$ws.Range("L122").Value = 15935.625
$ws.Range("M122").Value = -5678.3638
$ws.Range("N122").Value = -20835.625
$ws.Range("H126").Value = 2914.2856
$ws.Range("I126").Value = 1550
$ws.Range("J126").Value = 3460
$ws.Range("K126").Value = 4650
$ws.Range("L126").Value = 10380
$ws.Range("M126").Value = -2180
$ws.Range("N126").Value = -15320
$ws.Range("H136").Value = 3520.9111
$ws.Range("I136").Value = 2686.743
$ws.Range("J136").Value = 6440.5
$ws.Range("K136").Value = 8060.228999999999
$ws.Range("L136").Value = 19321.5
$ws.Range("M136").Value = -5510.228999999999
$ws.Range("N136").Value = -24421.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 41235.77
$ws.Range("I132").Value = 13882.75
$ws.Range("J132").Value = 85000.60000000001
$ws.Range("K132").Value = 41648.25
$ws.Range("L132").Value = 255001.8
$ws.Range("M132").Value = -39118.25
$ws.Range("N132").Value = -260061.8
$ws.Range("H136").Value = 2622.054
$ws.Range("I136").Value = 2024.0769
$ws.Range("J136").Value = 4035.4546
$ws.Range("K136").Value = 6072.2307
$ws.Range("L136").Value = 12106.3638
$ws.Range("M136").Value = -3522.2307
$ws.Range("N136").Value = -17206.3638

